$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Sema4a"
$ws.Cells.Item(2, 3).Value = "Plxnd1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 6.686327333333334
$ws.Cells.Item(2, 8).Value = 20.058982
$ws.Cells.Item(2, 9).Value = 0.1764498904644473
$ws.Cells.Item(2, 10).Value = 0.1764498904644472
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 81.979392
$ws.Cells.Item(2, 14).Value = 245.938176
$ws.Cells.Item(2, 15).Value = 0.523851352180617
$ws.Cells.Item(2, 16).Value = 0.523851352180617
$ws.Cells.Item(2, 17).Value = 548.141049499648
$ws.Cells.Item(2, 18).Value = 4933.269445496832
$ws.Cells.Item(2, 19).Value = 0.09243351371192245
$ws.Cells.Item(2, 20).Value = 0.09243351371192243

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Sema4a"
$ws.Cells.Item(3, 3).Value = "Plxnd1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 6.686327333333334
$ws.Cells.Item(3, 8).Value = 20.058982
$ws.Cells.Item(3, 9).Value = 0.1764498904644473
$ws.Cells.Item(3, 10).Value = 0.1764498904644472
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 22.17197066666667
$ws.Cells.Item(3, 14).Value = 66.515912
$ws.Cells.Item(3, 15).Value = 0.1416797140218155
$ws.Cells.Item(3, 16).Value = 0.1416797140218155
$ws.Cells.Item(3, 17).Value = 148.2490535023982
$ws.Cells.Item(3, 18).Value = 1334.241481521584
$ws.Cells.Item(3, 19).Value = 0.02499937002018355
$ws.Cells.Item(3, 20).Value = 0.02499937002018355

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Sema4a"
$ws.Cells.Item(4, 3).Value = "Plxnd1"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 6.686327333333334
$ws.Cells.Item(4, 8).Value = 20.058982
$ws.Cells.Item(4, 9).Value = 0.1764498904644473
$ws.Cells.Item(4, 10).Value = 0.1764498904644472
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 29.98794266666667
$ws.Cells.Item(4, 14).Value = 89.963828
$ws.Cells.Item(4, 15).Value = 0.1916240646801595
$ws.Cells.Item(4, 16).Value = 0.1916240646801595
$ws.Cells.Item(4, 17).Value = 200.5092007225662
$ws.Cells.Item(4, 18).Value = 1804.582806503096
$ws.Cells.Item(4, 19).Value = 0.0338120452231663
$ws.Cells.Item(4, 20).Value = 0.0338120452231663

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Sema4a"
$ws.Cells.Item(5, 3).Value = "Plxnd1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 6.686327333333334
$ws.Cells.Item(5, 8).Value = 20.058982
$ws.Cells.Item(5, 9).Value = 0.1764498904644473
$ws.Cells.Item(5, 10).Value = 0.1764498904644472
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 22.35431
$ws.Cells.Item(5, 14).Value = 67.06293000000001
$ws.Cells.Item(5, 15).Value = 0.142844869117408
$ws.Cells.Item(5, 16).Value = 0.1428448691174081
$ws.Cells.Item(5, 17).Value = 149.4682339708067
$ws.Cells.Item(5, 18).Value = 1345.21410573726
$ws.Cells.Item(5, 19).Value = 0.02520496150917495
$ws.Cells.Item(5, 20).Value = 0.02520496150917495

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Sema4a"
$ws.Cells.Item(6, 3).Value = "Plxnd1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 13.37868633333333
$ws.Cells.Item(6, 8).Value = 40.136059
$ws.Cells.Item(6, 9).Value = 0.353058954548371
$ws.Cells.Item(6, 10).Value = 0.353058954548371
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 81.979392
$ws.Cells.Item(6, 14).Value = 245.938176
$ws.Cells.Item(6, 15).Value = 0.523851352180617
$ws.Cells.Item(6, 16).Value = 0.523851352180617
$ws.Cells.Item(6, 17).Value = 1096.776571365376
$ws.Cells.Item(6, 18).Value = 9870.989142288385
$ws.Cells.Item(6, 19).Value = 0.1849504107396391
$ws.Cells.Item(6, 20).Value = 0.1849504107396391

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Sema4a"
$ws.Cells.Item(7, 3).Value = "Plxnd1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 13.37868633333333
$ws.Cells.Item(7, 8).Value = 40.136059
$ws.Cells.Item(7, 9).Value = 0.353058954548371
$ws.Cells.Item(7, 10).Value = 0.353058954548371
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 22.17197066666667
$ws.Cells.Item(7, 14).Value = 66.515912
$ws.Cells.Item(7, 15).Value = 0.1416797140218155
$ws.Cells.Item(7, 16).Value = 0.1416797140218155
$ws.Cells.Item(7, 17).Value = 296.6318409412009
$ws.Cells.Item(7, 18).Value = 2669.686568470808
$ws.Cells.Item(7, 19).Value = 0.05002129171325435
$ws.Cells.Item(7, 20).Value = 0.05002129171325436

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Sema4a"
$ws.Cells.Item(8, 3).Value = "Plxnd1"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 13.37868633333333
$ws.Cells.Item(8, 8).Value = 40.136059
$ws.Cells.Item(8, 9).Value = 0.353058954548371
$ws.Cells.Item(8, 10).Value = 0.353058954548371
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 29.98794266666667
$ws.Cells.Item(8, 14).Value = 89.963828
$ws.Cells.Item(8, 15).Value = 0.1916240646801595
$ws.Cells.Item(8, 16).Value = 0.1916240646801595
$ws.Cells.Item(8, 17).Value = 401.199278719317
$ws.Cells.Item(8, 18).Value = 3610.793508473852
$ws.Cells.Item(8, 19).Value = 0.06765459194228654
$ws.Cells.Item(8, 20).Value = 0.06765459194228654

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Sema4a"
$ws.Cells.Item(9, 3).Value = "Plxnd1"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 13.37868633333333
$ws.Cells.Item(9, 8).Value = 40.136059
$ws.Cells.Item(9, 9).Value = 0.353058954548371
$ws.Cells.Item(9, 10).Value = 0.353058954548371
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 22.35431
$ws.Cells.Item(9, 14).Value = 67.06293000000001
$ws.Cells.Item(9, 15).Value = 0.142844869117408
$ws.Cells.Item(9, 16).Value = 0.1428448691174081
$ws.Cells.Item(9, 17).Value = 299.0713016880967
$ws.Cells.Item(9, 18).Value = 2691.64171519287
$ws.Cells.Item(9, 19).Value = 0.05043266015319097
$ws.Cells.Item(9, 20).Value = 0.05043266015319098

$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Sema4a"
$ws.Cells.Item(10, 3).Value = "Plxnd1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 13.62475033333333
$ws.Cells.Item(10, 8).Value = 40.874251
$ws.Cells.Item(10, 9).Value = 0.3595524993125934
$ws.Cells.Item(10, 10).Value = 0.3595524993125934
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 81.979392
$ws.Cells.Item(10, 14).Value = 245.938176
$ws.Cells.Item(10, 15).Value = 0.523851352180617
$ws.Cells.Item(10, 16).Value = 0.523851352180617
$ws.Cells.Item(10, 17).Value = 1116.948748478464
$ws.Cells.Item(10, 18).Value = 10052.53873630618
$ws.Cells.Item(10, 19).Value = 0.1883520629448224
$ws.Cells.Item(10, 20).Value = 0.1883520629448224

$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Sema4a"
$ws.Cells.Item(11, 3).Value = "Plxnd1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 13.62475033333333
$ws.Cells.Item(11, 8).Value = 40.874251
$ws.Cells.Item(11, 9).Value = 0.3595524993125934
$ws.Cells.Item(11, 10).Value = 0.3595524993125934
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 22.17197066666667
$ws.Cells.Item(11, 14).Value = 66.515912
$ws.Cells.Item(11, 15).Value = 0.1416797140218155
$ws.Cells.Item(11, 16).Value = 0.1416797140218155
$ws.Cells.Item(11, 17).Value = 302.0875647313235
$ws.Cells.Item(11, 18).Value = 2718.788082581912
$ws.Cells.Item(11, 19).Value = 0.05094129527843724
$ws.Cells.Item(11, 20).Value = 0.05094129527843724

$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Sema4a"
$ws.Cells.Item(12, 3).Value = "Plxnd1"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 13.62475033333333
$ws.Cells.Item(12, 8).Value = 40.874251
$ws.Cells.Item(12, 9).Value = 0.3595524993125934
$ws.Cells.Item(12, 10).Value = 0.3595524993125934
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 29.98794266666667
$ws.Cells.Item(12, 14).Value = 89.963828
$ws.Cells.Item(12, 15).Value = 0.1916240646801595
$ws.Cells.Item(12, 16).Value = 0.1916240646801595
$ws.Cells.Item(12, 17).Value = 408.5782318436476
$ws.Cells.Item(12, 18).Value = 3677.204086592828
$ws.Cells.Item(12, 19).Value = 0.0688989113841894
$ws.Cells.Item(12, 20).Value = 0.06889891138418941

$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Sema4a"
$ws.Cells.Item(13, 3).Value = "Plxnd1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 13.62475033333333
$ws.Cells.Item(13, 8).Value = 40.874251
$ws.Cells.Item(13, 9).Value = 0.3595524993125934
$ws.Cells.Item(13, 10).Value = 0.3595524993125934
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 22.35431
$ws.Cells.Item(13, 14).Value = 67.06293000000001
$ws.Cells.Item(13, 15).Value = 0.142844869117408
$ws.Cells.Item(13, 16).Value = 0.1428448691174081
$ws.Cells.Item(13, 17).Value = 304.5718926239367
$ws.Cells.Item(13, 18).Value = 2741.14703361543
$ws.Cells.Item(13, 19).Value = 0.05136022970514435
$ws.Cells.Item(13, 20).Value = 0.05136022970514435

$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Sema4a"
$ws.Cells.Item(14, 3).Value = "Plxnd1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 4.203868666666668
$ws.Cells.Item(14, 8).Value = 12.611606
$ws.Cells.Item(14, 9).Value = 0.1109386556745884
$ws.Cells.Item(14, 10).Value = 0.1109386556745884
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 81.979392
$ws.Cells.Item(14, 14).Value = 245.938176
$ws.Cells.Item(14, 15).Value = 0.523851352180617
$ws.Cells.Item(14, 16).Value = 0.523851352180617
$ws.Cells.Item(14, 17).Value = 344.6305973411841
$ws.Cells.Item(14, 18).Value = 3101.675376070657
$ws.Cells.Item(14, 19).Value = 0.05811536478423299
$ws.Cells.Item(14, 20).Value = 0.05811536478423299

$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Sema4a"
$ws.Cells.Item(15, 3).Value = "Plxnd1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 4.203868666666668
$ws.Cells.Item(15, 8).Value = 12.611606
$ws.Cells.Item(15, 9).Value = 0.1109386556745884
$ws.Cells.Item(15, 10).Value = 0.1109386556745884
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 22.17197066666667
$ws.Cells.Item(15, 14).Value = 66.515912
$ws.Cells.Item(15, 15).Value = 0.1416797140218155
$ws.Cells.Item(15, 16).Value = 0.1416797140218155
$ws.Cells.Item(15, 17).Value = 93.20805276385246
$ws.Cells.Item(15, 18).Value = 838.8724748746721
$ws.Cells.Item(15, 19).Value = 0.01571775700994034
$ws.Cells.Item(15, 20).Value = 0.01571775700994034

$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Sema4a"
$ws.Cells.Item(16, 3).Value = "Plxnd1"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 4.203868666666668
$ws.Cells.Item(16, 8).Value = 12.611606
$ws.Cells.Item(16, 9).Value = 0.1109386556745884
$ws.Cells.Item(16, 10).Value = 0.1109386556745884
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 29.98794266666667
$ws.Cells.Item(16, 14).Value = 89.963828
$ws.Cells.Item(16, 15).Value = 0.1916240646801595
$ws.Cells.Item(16, 16).Value = 0.1916240646801595
$ws.Cells.Item(16, 17).Value = 126.0653725541965
$ws.Cells.Item(16, 18).Value = 1134.588352987768
$ws.Cells.Item(16, 19).Value = 0.02125851613051727
$ws.Cells.Item(16, 20).Value = 0.02125851613051727

$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Sema4a"
$ws.Cells.Item(17, 3).Value = "Plxnd1"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 4.203868666666668
$ws.Cells.Item(17, 8).Value = 12.611606
$ws.Cells.Item(17, 9).Value = 0.1109386556745884
$ws.Cells.Item(17, 10).Value = 0.1109386556745884
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 22.35431
$ws.Cells.Item(17, 14).Value = 67.06293000000001
$ws.Cells.Item(17, 15).Value = 0.142844869117408
$ws.Cells.Item(17, 16).Value = 0.1428448691174081
$ws.Cells.Item(17, 17).Value = 93.97458337395337
$ws.Cells.Item(17, 18).Value = 845.7712503655803
$ws.Cells.Item(17, 19).Value = 0.01584701774989777
$ws.Cells.Item(17, 20).Value = 0.01584701774989777
